$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 210, pushing the existing rows 210..315 down to 212..317
$ws.Rows.Item(210).Resize(2).Insert()

# New weekly records, newly inserted at rows 210 and 211
$newRows = @(
    @{ Row = 210; A = 10; B = "Vega Modelo de Temuco"; C = "La Araucanía"; D = 45001; E = 9; F = 100112043; G = "Pepino dulce"; H = "Cultivar XV región"; I = "Extra";   J = 100; K = 17000; L = 17000; M = 17000; N = "`$/bandeja 18 kilos"; O = "Región de Arica y Parinacota"; P = 944; Q = 18; R = "Hortaliza" },
    @{ Row = 211; A = 10; B = "Vega Modelo de Temuco"; C = "La Araucanía"; D = 45001; E = 9; F = 100112043; G = "Pepino dulce"; H = "Cultivar XV región"; I = "Primera";  J = 500; K = 14000; L = 15000; M = 14600; N = "`$/bandeja 18 kilos"; O = "Región de Arica y Parinacota"; P = 811; Q = 18; R = "Hortaliza" }
)

foreach ($rec in $newRows) {
    $r = $rec.Row
    $ws.Cells.Item($r, 1).Value = $rec.A
    $ws.Cells.Item($r, 2).Value = $rec.B
    $ws.Cells.Item($r, 3).Value = $rec.C
    $ws.Cells.Item($r, 4).Value = $rec.D
    $ws.Cells.Item($r, 5).Value = $rec.E
    $ws.Cells.Item($r, 6).Value = $rec.F
    $ws.Cells.Item($r, 7).Value = $rec.G
    $ws.Cells.Item($r, 8).Value = $rec.H
    $ws.Cells.Item($r, 9).Value = $rec.I
    $ws.Cells.Item($r, 10).Value = $rec.J
    $ws.Cells.Item($r, 11).Value = $rec.K
    $ws.Cells.Item($r, 12).Value = $rec.L
    $ws.Cells.Item($r, 13).Value = $rec.M
    $ws.Cells.Item($r, 14).Value = $rec.N
    $ws.Cells.Item($r, 15).Value = $rec.O
    $ws.Cells.Item($r, 16).Value = $rec.P
    $ws.Cells.Item($r, 17).Value = $rec.Q
    $ws.Cells.Item($r, 18).Value = $rec.R
}
